# Exercise 8 ("Realizacao do exercicio 8"): append the JSHELL transcript
# that demonstrates running the "preco do carro" pseudocode algorithm.
#
# The transcript is inserted right after the LAST "No JSHELL" paragraph
# in the document (the one that closes exercise 8), before the blank
# paragraph that precedes the next exercise's list item.

$d = $word.ActiveDocument

# Locate the paragraph that holds the final "No JSHELL" marker (exercise 8).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "No JSHELL") {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'No JSHELL' paragraph to anchor the new content on."
}

$anchor = $d.Paragraphs.Item($targetIndex).Range
$anchor.Collapse(0)

# Build the JSHELL session transcript as one block of text, using carriage
# returns to delimit the individual paragraphs that need to be created.
# (Word COM turns each embedded vbCr into its own new paragraph, inheriting
# the paragraph formatting -- spacing after = 0 -- from the anchor paragraph.)
$lines = @(
    "",
    'jshell> System.out.println("Favor informar o modelo do carro");',
    "Favor informar o modelo do carro",
    "",
    "jshell> String modelo = teclado.next();",
    "Corolla",
    'titulo_filme ==> "Corolla"',
    "",
    "jshell> double preco_carro = teclado.nextDouble();",
    "80000",
    "preco_carro ==> 80000.0",
    "",
    "jshell> double entrada = preco_carro80.5;",
    "entrada ==> 40000.0",
    "",
    "jshell> double restante = preco_carro*0.5/12;",
    "restante ==> 3333.3333333333335",
    "",
    "jshell>"
)

$block = [string]::Join([char]13, $lines)

$anchor.InsertAfter($block)
